$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values per row for columns D, L, M, N, O, P, Q, R, S, T
# (derived from a reordering of the weekly price records)
$data = @{
    3 = @{ D=44515; L='Primera'; M=80; N=28000; O=28000; P=28000; Q='$/bandeja 10 kilos'; R='Provincia de Los Andes'; S=2800; T=10 }
    4 = @{ D=44859; L='Primera'; M=30; N=20000; O=20000; P=20000; Q='$/bandeja 5 kilos'; R='Provincia de Quillota'; S=4000; T=5 }
    5 = @{ D=44496; L='Primera'; M=55; N=28000; O=28000; P=28000; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=2800; T=10 }
    6 = @{ D=44466; L='Primera'; M=80; N=11000; O=11000; P=11000; Q='$/bandeja 5 kilos'; R='La Ligua'; S=2200; T=5 }
    7 = @{ D=44511; L='Primera'; M=45; N=28000; O=28000; P=28000; Q='$/bandeja 10 kilos'; R='Provincia de Los Andes'; S=2800; T=10 }
    8 = @{ D=44511; L='Primera'; M=45; N=3200; O=3200; P=3200; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=320; T=10 }
    9 = @{ D=44874; L='Primera'; M=40; N=25000; O=25000; P=25000; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=2500; T=10 }
    10 = @{ D=44902; L='Primera'; M=90; N=25000; O=25000; P=25000; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=2500; T=10 }
    11 = @{ D=44166; L='Segunda'; M=20; N=12000; O=12000; P=12000; Q='$/caja 18 kilos'; R='La Ligua'; S=667; T=18 }
    12 = @{ D=44503; L='Primera'; M=50; N=28000; O=28000; P=28000; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=2800; T=10 }
    13 = @{ D=44921; L='Primera'; M=55; N=15000; O=15000; P=15000; Q='$/bandeja 7 kilos'; R='Provincia de Quillota'; S=2143; T=7 }
    14 = @{ D=44901; L='Primera'; M=40; N=25000; O=25000; P=25000; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=2500; T=10 }
    15 = @{ D=44889; L='Primera'; M=50; N=30000; O=30000; P=30000; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=3000; T=10 }
    16 = @{ D=44519; L='Primera'; M=30; N=28000; O=28000; P=28000; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=2800; T=10 }
    17 = @{ D=44858; L='Primera'; M=90; N=20000; O=20000; P=20000; Q='$/bandeja 5 kilos'; R='Provincia de Quillota'; S=4000; T=5 }
    18 = @{ D=44868; L='Primera'; M=30; N=14000; O=14000; P=14000; Q='$/bandeja 5 kilos'; R='Provincia de Quillota'; S=2800; T=5 }
    19 = @{ D=44879; L='Primera'; M=25; N=30000; O=30000; P=30000; Q='$/bandeja 10 kilos'; R='Provincia de Quillota'; S=3000; T=10 }
    20 = @{ D=44488; L='Primera'; M=100; N=12000; O=12000; P=12000; Q='$/bandeja 5 kilos'; R='La Ligua'; S=2400; T=5 }
}

foreach ($r in $data.Keys) {
    $row = [int]$r
    $vals = $data[$r]
    $ws.Cells.Item($row, 4).Value  = $vals.D   # D: Fecha
    $ws.Cells.Item($row, 12).Value = $vals.L   # L: Calidad
    $ws.Cells.Item($row, 13).Value = $vals.M   # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N   # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O   # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P   # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals.Q   # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $vals.R   # R: Origen
    $ws.Cells.Item($row, 19).Value = $vals.S   # S: Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $vals.T   # T: Kg / unidad
}